# Update imputed values produced by the RandomForest algorithm run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.4409
$ws.Range("C3").Value = -11.5951
$ws.Range("E4").Value = 12.88630000000001
$ws.Range("C5").Value = -12.63740000000001
$ws.Range("E6").Value = 11.99549999999999
$ws.Range("D7").Value = -7.15339999999999
$ws.Range("B9").Value = 8.469100000000006
$ws.Range("D9").Value = -8.754100000000001
$ws.Range("E10").Value = 11.49009999999999
$ws.Range("C11").Value = -13.4717
$ws.Range("E11").Value = 13.3879
$ws.Range("C12").Value = -14.31590000000002
$ws.Range("B13").Value = 5.036600000000004
$ws.Range("B16").Value = 9.142800000000008
$ws.Range("B18").Value = 4.841700000000004
$ws.Range("B20").Value = 5.599899999999995
$ws.Range("C21").Value = -13.0609
$ws.Range("D21").Value = -8.181400000000004
$ws.Range("E21").Value = 13.0726
$ws.Range("E25").Value = 13.4309
